# Update scripts with new TPM values.
# Adds a new 'MuSCs' sending-cluster block (rows 7-11) mirroring the existing
# Gdf1->Acvr1 interactions, and refreshes the recomputed specificity/weight
# values on the original FAPs rows (2-6) to reflect the new TPM normalization.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Gdf1"
$ws.Cells.Item(2,3).Value = "Acvr1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 0.1814253333333333
$ws.Cells.Item(2,8).Value = 0.544276
$ws.Cells.Item(2,9).Value = 0.9591577789839493
$ws.Cells.Item(2,10).Value = 0.9591577789839494
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 9.358310000000001
$ws.Cells.Item(2,14).Value = 28.07493
$ws.Cells.Item(2,15).Value = 0.1719944618809179
$ws.Cells.Item(2,16).Value = 0.1719944618809179
$ws.Cells.Item(2,17).Value = 1.697834511186667
$ws.Cells.Item(2,18).Value = 15.28051060068
$ws.Cells.Item(2,19).Value = 0.1649698260552407
$ws.Cells.Item(2,20).Value = 0.1649698260552408

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Gdf1"
$ws.Cells.Item(3,3).Value = "Acvr1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 0.1814253333333333
$ws.Cells.Item(3,8).Value = 0.544276
$ws.Cells.Item(3,9).Value = 0.9591577789839493
$ws.Cells.Item(3,10).Value = 0.9591577789839494
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 26.47935433333333
$ws.Cells.Item(3,14).Value = 79.438063
$ws.Cells.Item(3,15).Value = 0.4866586274141183
$ws.Cells.Item(3,16).Value = 0.4866586274141184
$ws.Cells.Item(3,17).Value = 4.804025686376444
$ws.Cells.Item(3,18).Value = 43.236231177388
$ws.Cells.Item(3,19).Value = 0.466782408193903
$ws.Cells.Item(3,20).Value = 0.4667824081939031

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Gdf1"
$ws.Cells.Item(4,3).Value = "Acvr1"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 0.1814253333333333
$ws.Cells.Item(4,8).Value = 0.544276
$ws.Cells.Item(4,9).Value = 0.9591577789839493
$ws.Cells.Item(4,10).Value = 0.9591577789839494
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.152806666666667
$ws.Cells.Item(4,14).Value = 15.45842
$ws.Cells.Item(4,15).Value = 0.09470237786627494
$ws.Cells.Item(4,16).Value = 0.09470237786627496
$ws.Cells.Item(4,17).Value = 0.9348496671022222
$ws.Cells.Item(4,18).Value = 8.41364700392
$ws.Cells.Item(4,19).Value = 0.090834522418715
$ws.Cells.Item(4,20).Value = 0.09083452241871502

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Gdf1"
$ws.Cells.Item(5,3).Value = "Acvr1"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.1814253333333333
$ws.Cells.Item(5,8).Value = 0.544276
$ws.Cells.Item(5,9).Value = 0.9591577789839493
$ws.Cells.Item(5,10).Value = 0.9591577789839494
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 7.974813333333334
$ws.Cells.Item(5,14).Value = 23.92444
$ws.Cells.Item(5,15).Value = 0.146567460136225
$ws.Cells.Item(5,16).Value = 0.146567460136225
$ws.Cells.Item(5,17).Value = 1.446833167271111
$ws.Cells.Item(5,18).Value = 13.02149850544
$ws.Cells.Item(5,19).Value = 0.1405813195355801
$ws.Cells.Item(5,20).Value = 0.1405813195355801

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Gdf1"
$ws.Cells.Item(6,3).Value = "Acvr1"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.1814253333333333
$ws.Cells.Item(6,8).Value = 0.544276
$ws.Cells.Item(6,9).Value = 0.9591577789839493
$ws.Cells.Item(6,10).Value = 0.9591577789839494
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 5.445246666666667
$ws.Cells.Item(6,14).Value = 16.33574
$ws.Cells.Item(6,15).Value = 0.1000770727024639
$ws.Cells.Item(6,16).Value = 0.1000770727024639
$ws.Cells.Item(6,17).Value = 0.9879056915822222
$ws.Cells.Item(6,18).Value = 8.89115122424
$ws.Cells.Item(6,19).Value = 0.09598970278051051
$ws.Cells.Item(6,20).Value = 0.09598970278051054

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Gdf1"
$ws.Cells.Item(7,3).Value = "Acvr1"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.007725333333333334
$ws.Cells.Item(7,8).Value = 0.023176
$ws.Cells.Item(7,9).Value = 0.0408422210160507
$ws.Cells.Item(7,10).Value = 0.0408422210160507
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 9.358310000000001
$ws.Cells.Item(7,14).Value = 28.07493
$ws.Cells.Item(7,15).Value = 0.1719944618809179
$ws.Cells.Item(7,16).Value = 0.1719944618809179
$ws.Cells.Item(7,17).Value = 0.07229606418666669
$ws.Cells.Item(7,18).Value = 0.6506645776800001
$ws.Cells.Item(7,19).Value = 0.007024635825677156
$ws.Cells.Item(7,20).Value = 0.007024635825677156

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Gdf1"
$ws.Cells.Item(8,3).Value = "Acvr1"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.007725333333333334
$ws.Cells.Item(8,8).Value = 0.023176
$ws.Cells.Item(8,9).Value = 0.0408422210160507
$ws.Cells.Item(8,10).Value = 0.0408422210160507
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 26.47935433333333
$ws.Cells.Item(8,14).Value = 79.438063
$ws.Cells.Item(8,15).Value = 0.4866586274141183
$ws.Cells.Item(8,16).Value = 0.4866586274141184
$ws.Cells.Item(8,17).Value = 0.2045618386764445
$ws.Cells.Item(8,18).Value = 1.841056548088
$ws.Cells.Item(8,19).Value = 0.01987621922021529
$ws.Cells.Item(8,20).Value = 0.01987621922021529

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Gdf1"
$ws.Cells.Item(9,3).Value = "Acvr1"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.007725333333333334
$ws.Cells.Item(9,8).Value = 0.023176
$ws.Cells.Item(9,9).Value = 0.0408422210160507
$ws.Cells.Item(9,10).Value = 0.0408422210160507
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 5.152806666666667
$ws.Cells.Item(9,14).Value = 15.45842
$ws.Cells.Item(9,15).Value = 0.09470237786627494
$ws.Cells.Item(9,16).Value = 0.09470237786627496
$ws.Cells.Item(9,17).Value = 0.03980714910222222
$ws.Cells.Item(9,18).Value = 0.3582643419200001
$ws.Cells.Item(9,19).Value = 0.003867855447559949
$ws.Cells.Item(9,20).Value = 0.00386785544755995

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Gdf1"
$ws.Cells.Item(10,3).Value = "Acvr1"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.007725333333333334
$ws.Cells.Item(10,8).Value = 0.023176
$ws.Cells.Item(10,9).Value = 0.0408422210160507
$ws.Cells.Item(10,10).Value = 0.0408422210160507
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 7.974813333333334
$ws.Cells.Item(10,14).Value = 23.92444
$ws.Cells.Item(10,15).Value = 0.146567460136225
$ws.Cells.Item(10,16).Value = 0.146567460136225
$ws.Cells.Item(10,17).Value = 0.06160809127111112
$ws.Cells.Item(10,18).Value = 0.5544728214400001
$ws.Cells.Item(10,19).Value = 0.005986140600644901
$ws.Cells.Item(10,20).Value = 0.005986140600644903

# Row 11
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Gdf1"
$ws.Cells.Item(11,3).Value = "Acvr1"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.007725333333333334
$ws.Cells.Item(11,8).Value = 0.023176
$ws.Cells.Item(11,9).Value = 0.0408422210160507
$ws.Cells.Item(11,10).Value = 0.0408422210160507
$ws.Cells.Item(11,11).Value = 3
$ws.Cells.Item(11,12).Value = 1
$ws.Cells.Item(11,13).Value = 5.445246666666667
$ws.Cells.Item(11,14).Value = 16.33574
$ws.Cells.Item(11,15).Value = 0.1000770727024639
$ws.Cells.Item(11,16).Value = 0.1000770727024639
$ws.Cells.Item(11,17).Value = 0.04206634558222223
$ws.Cells.Item(11,18).Value = 0.3785971102400001
$ws.Cells.Item(11,19).Value = 0.004087369921953406
$ws.Cells.Item(11,20).Value = 0.004087369921953406
